# Add a "Save" column (H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: copy formatting (bold + border, style "s=1") from the
# neighbouring header cell G1, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Data cells H2:H32: the "Save" indicator values for each row.
$values = @(0,1,1,0,0,1,1,0,1,0,0,0,0,0,0,0,0,0,1,0,1,0,0,0,0,0,0,0,0,0,0)
$r = 2
foreach ($v in $values) {
    $ws.Cells.Item($r, 8).Value = $v
    $r = $r + 1
}
